$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.505.93"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "2.578.72"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.53"
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.69"
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -1.22%  "
$ws.Range("D9").Value = "2.577.38"
$ws.Range("E9").Value = "  -1.58%  "
$ws.Range("E10").Value = "  -3.24%  "
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.359"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.17"
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.79"
$ws.Range("E14").Value = "  -1.96%  "
$ws.Range("D15").Value = "3.048.11"
$ws.Range("E15").Value = "  -2.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000178"
$ws.Range("E16").Value = "  -1.56%  "
$ws.Range("D17").Value = "66.648.09"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").Value = "2.594.41"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("E19").Value = "  -4.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.74"
$ws.Range("E20").Value = "  -3.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "351.65"
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("E22").Value = "  -2.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.60"
$ws.Range("E23").Value = "  -1.38%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  -2.98%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.02"
$ws.Range("E26").Value = "  -2.24%  "
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.02"
$ws.Range("E27").Value = "  -7.85%  "
$ws.Range("D28").Value = "2.710.55"
$ws.Range("E28").Value = "  -1.63%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "0.0$([char]0x2083)0989"
$ws.Range("E30").Value = "  -1.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "536.16"
$ws.Range("E31").Value = "  -2.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.03"
$ws.Range("E32").Value = "  +1.75%  "
$ws.Range("E33").Value = "  -1.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.85"
$ws.Range("E34").Value = "  -2.33%  "
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  -3.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.84"
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.75"
$ws.Range("E39").Value = "  -1.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.361"
$ws.Range("E40").Value = "  -1.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.24"
$ws.Range("E41").Value = "  +1.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.79"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.11"
$ws.Range("E43").Value = "  -1.40%  "
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.41"
$ws.Range("E45").Value = "  -1.86%  "
$ws.Range("D46").Value = "0.0$([char]0x2086)0288"
$ws.Range("E46").Value = "  -2.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "148.93"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.567"
$ws.Range("E48").Value = "  -2.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.73"
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("E51").Value = "  -1.72%  "
